# Update the "Fund Source Data" sheet's D1 header from
# "Distribution Amount" to "Distribution" to align with the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fund Source Data")
$ws.Range("D1").Value = "Distribution"
